$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J column (k values)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary stats
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold, larger font for the summary values, vertically centered
$range = $ws.Range("B14:B17")
$range.Font.Bold = $true
$range.Font.Size = 12
$range.VerticalAlignment = -4108

# Selection / print setup to mirror the final interactive state
$ws.Range("A14:B17").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
